$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing giraffe entry's accession-ID typo (Erthryo -> Erythyro) ---
# Doing this first, while the old text is still referenced only by A9, lets
# the corrected string take over the same shared-string slot the old one used.
$ws.Range("A9").Value = "ErythyroPV-Giraffe-1"

# --- New row 10: novel hyaena erythyroparvovirus identified from WGS mining ---
$ws.Range("A10").Value = "ErythyroPV-Hyaena-1"

# --- Host species, now known for these two Ungulate erythroparvovirus-2 records ---
$ws.Range("D9").Value  = "Giraffa camelopardalis"
$ws.Range("D10").Value = "Hyaena hyaena"

# --- Host species, now known for the existing human B19V record ---
$ws.Range("D3").Value = "Homo sapiens"

$ws.Range("B10").Value = "HhEPV"

$ws.Range("C10").Value = "Ungulate erythroparvovirus 2"
$ws.Range("E10").Value = "NK"
$ws.Range("F10").Value = "Erythyroparvovirus"
$ws.Range("G10").Value = "NK"
$ws.Range("H10").Value = "NK"
$ws.Range("I10").Value = "NK"
$ws.Range("J10").Value = "NK"
$ws.Range("K10").Value = "NK"
$ws.Range("L10").Value = $false

# --- Widen column A to fit the longer accession-ID labels ---
$ws.Range("A1").ColumnWidth = 24

# --- Restore the sheet's selection state (whole table selected) ---
[void]$ws.Range("A1:L10").Select()
[void]$ws.Range("C7").Activate()
